$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.0425
$ws.Range("E2").Value = -0.0458
$ws.Range("G2").Value = 0.139544977582753
$ws.Range("H2").Value = 0.139544977582753
$ws.Range("I2").Value = 0.04062291328818086
$ws.Range("J2").Value = 0.0296642755187035
$ws.Range("K2").Value = 59.66
$ws.Range("L2").Value = 0.02845559477248879
$ws.Range("M2").Value = 81.2
$ws.Range("N2").Value = 0.02825232246616332
$ws.Range("O2").Value = 1.361045926919209
$ws.Range("P2").Value = 81.2
$ws.Range("Q2").Value = 0.02825232246616332
$ws.Range("R2").Value = 1.361045926919209
$ws.Range("U2").Value = 61.125
$ws.Range("V2").Value = 0.0212675272259142
$ws.Range("W2").Value = 0.2372372372372372
$ws.Range("X2").Value = 0.05774411899772042
$ws.Range("Y2").Value = 0.1794931182395168
$ws.Range("Z2").Value = 1.884016947705635
$ws.Range("AA2").Value = 0.01081730769230769
$ws.Range("AB2").Value = 0.05666118767301267
$ws.Range("AC2").Value = -0.04584387998070498
$ws.Range("AD2").Value = 94.484
$ws.Range("AF2").Value = 94.484
$ws.Range("AG2").Value = 33.35899999999999
$ws.Range("AH2").Value = 0.03182796915970712
$ws.Range("AI2").Value = 0.07241129838200401
$ws.Range("AJ2").Value = 0.01147359257688586
$ws.Range("AK2").Value = 0.02682240638611111
$ws.Range("AL2").Value = 0.547
$ws.Range("AM2").Value = 0.294
$ws.Range("AN2").Value = 0.3964419082784374
$ws.Range("AO2").Value = 155.7038391224863
$ws.Range("AP2").Value = 0.1399697897872697
$ws.Range("AQ2").Value = 289.6938775510204
$ws.Range("D3").Value = 0.0262
$ws.Range("E3").Value = -0.0458
$ws.Range("G3").Value = 0.4346534653465347
$ws.Range("H3").Value = 0.4346534653465347
$ws.Range("I3").Value = 0.4222772277227723
$ws.Range("J3").Value = 0.2913100875304921
$ws.Range("K3").Value = 4.76
$ws.Range("L3").Value = 0.2356435643564356
$ws.Range("U3").Value = 1.94
$ws.Range("V3").Value = 0.0154828411811652
$ws.Range("W3").Value = 1.919354838709677
$ws.Range("X3").Value = 0.06000091827922616
$ws.Range("Y3").Value = 1.859353920430451
$ws.Range("Z3").Value = 4.225941422594142
$ws.Range("AA3").Value = 1.231059365714632
$ws.Range("AB3").Value = 0.05693569586009507
$ws.Range("AC3").Value = 1.174123669854537
$ws.Range("AD3").Value = 13.5
$ws.Range("AF3").Value = 13.5
$ws.Range("AG3").Value = 11.56
$ws.Range("AH3").Value = 0.0972622478386167
$ws.Range("AI3").Value = 0.8743523316062176
$ws.Range("AJ3").Value = 0.08446587753909106
$ws.Range("AK3").Value = 0.8562962962962963
$ws.Range("AL3").Value = 0.547
$ws.Range("AM3").Value = 0.547
$ws.Range("AN3").Value = 1.527149321266968
$ws.Range("AO3").Value = 15.59414990859232
$ws.Range("AP3").Value = 1.307692307692308
$ws.Range("AQ3").Value = 15.59414990859232
$ws.Range("B4").Value = 'AtlantaSanad Société Anonyme (CBSE:ATL)'
$ws.Range("D4").Value = 0.09119999999999999
$ws.Range("E4").Value = 0.157
$ws.Range("G4").Value = 0.1379690949227373
$ws.Range("H4").Value = 0.1379690949227373
$ws.Range("I4").Value = 0.1035688005886681
$ws.Range("J4").Value = 0.07604029039502586
$ws.Range("K4").Value = 39.5
$ws.Range("L4").Value = 0.07266372332597498
$ws.Range("M4").Value = 20.5
$ws.Range("N4").Value = 0.03699693196173975
$ws.Range("O4").Value = 0.5189873417721519
$ws.Range("P4").Value = 20.5
$ws.Range("Q4").Value = 0.03699693196173975
$ws.Range("R4").Value = 0.5189873417721519
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 5.44
$ws.Range("V4").Value = 0.009817722432773868
$ws.Range("W4").Value = 0.2372372372372372
$ws.Range("X4").Value = 0.05989516561036436
$ws.Range("Y4").Value = 0.1773420716268729
$ws.Range("Z4").Value = 2.728641702640297
$ws.Range("AA4").Value = 0.207486707452746
$ws.Range("AB4").Value = 0.05692360186040313
$ws.Range("AC4").Value = 0.1505631055923428
$ws.Range("AD4").Value = 57.9
$ws.Range("AF4").Value = 57.9
$ws.Range("AG4").Value = 52.46
$ws.Range("AH4").Value = 0.0946078431372549
$ws.Range("AI4").Value = 0.239454094292804
$ws.Range("AJ4").Value = 0.08648773410709575
$ws.Range("AK4").Value = 0.2219495684548993
$ws.Range("AM4").Value = 0
$ws.Range("AN4").Value = 0.7470967741935484
$ws.Range("AP4").Value = 0.6769032258064516
$ws.Range("B5").Value = 'SAHAM Assurance S.A. (CBSE:SAH)'
$ws.Range("D5").Value = 0.0425
$ws.Range("E5").Value = -0.278
$ws.Range("G5").Value = 0.1640899508081518
$ws.Range("H5").Value = 0.1640899508081518
$ws.Range("I5").Value = 0.006324666198172874
$ws.Range("J5").Value = 0.006324666198172874
$ws.Range("K5").Value = 7.59
$ws.Range("L5").Value = 0.01333450456781447
$ws.Range("M5").Value = 12.7
$ws.Range("N5").Value = 0.02150355570606163
$ws.Range("O5").Value = 1.673254281949934
$ws.Range("P5").Value = 12.7
$ws.Range("Q5").Value = 0.02150355570606163
$ws.Range("R5").Value = 1.673254281949934
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 19.5
$ws.Range("V5").Value = 0.03301727057229936
$ws.Range("W5").Value = 0.01662650602409638
$ws.Range("X5").Value = 0.05774411899772042
$ws.Range("Y5").Value = -0.04111761297362403
$ws.Range("Z5").Value = 1.710336538461538
$ws.Range("AA5").Value = 0.01081730769230769
$ws.Range("AB5").Value = 0.05666118767301267
$ws.Range("AC5").Value = -0.04584387998070498
$ws.Range("AD5").Value = 22.7
$ws.Range("AF5").Value = 22.7
$ws.Range("AG5").Value = 3.199999999999999
$ws.Range("AH5").Value = 0.03701288113484428
$ws.Range("AI5").Value = 0.04713455149501661
$ws.Range("AJ5").Value = 0.005389019872010776
$ws.Range("AK5").Value = 0.006924908028565244
$ws.Range("AN5").Value = 0.6037234042553191
$ws.Range("AP5").Value = 0.08510638297872339
$ws.Range("D6").Value = 0.0513
$ws.Range("E6").Value = -0.4970000000000001
$ws.Range("G6").Value = 0.1137655230477794
$ws.Range("H6").Value = 0.1137655230477794
$ws.Range("I6").Value = 0.01015575668280362
$ws.Range("J6").Value = 0.005483712672780902
$ws.Range("K6").Value = 2.78
$ws.Range("L6").Value = 0.002925699852662597
$ws.Range("M6").Value = 43.3
$ws.Range("N6").Value = 0.02847560173615678
$ws.Range("O6").Value = 15.57553956834532
$ws.Range("P6").Value = 43.3
$ws.Range("Q6").Value = 0.02847560173615678
$ws.Range("R6").Value = 15.57553956834532
$ws.Range("U6").Value = 33.3
$ws.Range("V6").Value = 0.02189925029593581
$ws.Range("W6").Value = 0.004604173567406426
$ws.Range("X6").Value = 0.05650077273071079
$ws.Range("Y6").Value = -0.05189659916330436
$ws.Range("Z6").Value = 1.633979622544173
$ws.Range("AA6").Value = 0.008960274763211235
$ws.Range("AB6").Value = 0.05649369986558906
$ws.Range("AC6").Value = -0.04753342510237783
$ws.Range("AD6").Value = 0.384
$ws.Range("AF6").Value = 0.384
$ws.Range("AG6").Value = -32.916
$ws.Range("AH6").Value = 0.0002524681390468276
$ws.Range("AI6").Value = 0.0006911646123718465
$ws.Range("AJ6").Value = -0.02212566647218092
$ws.Range("AK6").Value = -0.0630231827894402
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = -0.019
$ws.Range("AN6").Value = 0.003585434173669468
$ws.Range("AP6").Value = -0.3073389355742297
$ws.Range("AQ6").Value = -507.8947368421053
$ws.Range("B7").Value = 'Agma S.A. (CBSE:AGM)'
$ws.Range("D7").Value = 0.0252
$ws.Range("E7").Value = 0.0134
$ws.Range("G7").Value = 0.5440298507462686
$ws.Range("H7").Value = 0.5440298507462686
$ws.Range("I7").Value = 0.5291044776119402
$ws.Range("J7").Value = 0.3635786232770573
$ws.Range("K7").Value = 5.03
$ws.Range("L7").Value = 0.3753731343283582
$ws.Range("M7").Value = 4.7
$ws.Range("N7").Value = 0.0562874251497006
$ws.Range("O7").Value = 0.9343936381709741
$ws.Range("P7").Value = 4.7
$ws.Range("Q7").Value = 0.0562874251497006
$ws.Range("R7").Value = 0.9343936381709741
$ws.Range("U7").Value = 0.945
$ws.Range("V7").Value = 0.01131736526946108
$ws.Range("W7").Value = 0.5085945399393327
$ws.Range("X7").Value = 0.05649254957015962
$ws.Range("Y7").Value = 0.4521019903691731
$ws.Range("Z7").Value = -2.4408014571949
$ws.Range("AA7").Value = -0.8874232334995572
$ws.Range("AB7").Value = 0.05649254957015962
$ws.Range("AC7").Value = -0.9439157830697168
$ws.Range("AD7").Value = 0
$ws.Range("AF7").Value = 0
$ws.Range("AG7").Value = -0.945
$ws.Range("AH7").Value = 0
$ws.Range("AI7").Value = 0
$ws.Range("AJ7").Value = -0.0114469141784265
$ws.Range("AK7").Value = -0.09994711792702274
$ws.Range("AL7").Value = 0
$ws.Range("AM7").Value = -0.234
$ws.Range("AN7").Value = 0
$ws.Range("AP7").Value = -0.1296296296296296
$ws.Range("AQ7").Value = -30.2991452991453
$ws.Range("AQ4").ClearContents()
$ws.Range("AO6").ClearContents()
$ws.Range("AO7").ClearContents()
